# Auto-generated edit script: updates cached market-price / profit figures
# across the per-class Leve profit worksheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 602.4286
$ws.Range("I9").Value = 538.6
$ws.Range("K9").Value = 538.6
$ws.Range("M9").Value = -369.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2444759.2
$ws.Range("I32").Value = 2505203.2
$ws.Range("J32").Value = 26999.5
$ws.Range("K32").Value = 2505203.2
$ws.Range("L32").Value = 26999.5
$ws.Range("M32").Value = -2504916.2
$ws.Range("N32").Value = -27573.5
$ws.Range("H45").Value = 4613.524
$ws.Range("I45").Value = 2302.5
$ws.Range("K45").Value = 2302.5
$ws.Range("M45").Value = -1925.5
$ws.Range("H122").Value = 1960.64
$ws.Range("I122").Value = 1617.5714
$ws.Range("K122").Value = 4852.7142
$ws.Range("M122").Value = -2402.7142
$ws.Range("H132").Value = 3337.547
$ws.Range("I132").Value = 1628.0294
$ws.Range("K132").Value = 4884.0882
$ws.Range("M132").Value = -2354.0882

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6176686.5
$ws.Range("I20").Value = 7249449
$ws.Range("J20").Value = 8302.25
$ws.Range("K20").Value = 7249449
$ws.Range("L20").Value = 8302.25
$ws.Range("M20").Value = -7249202
$ws.Range("N20").Value = -8796.25
$ws.Range("H86").Value = 49060.184
$ws.Range("I86").Value = 69702.2
$ws.Range("J86").Value = 4827.2856
$ws.Range("K86").Value = 69702.2
$ws.Range("L86").Value = 4827.2856
$ws.Range("M86").Value = -68579.2
$ws.Range("N86").Value = -7073.2856
$ws.Range("H89").Value = 49060.184
$ws.Range("I89").Value = 69702.2
$ws.Range("J89").Value = 4827.2856
$ws.Range("K89").Value = 348511
$ws.Range("L89").Value = 24136.428
$ws.Range("M89").Value = -342895
$ws.Range("N89").Value = -35368.428
$ws.Range("H113").Value = 4965.1665
$ws.Range("I113").Value = 4965.1665
$ws.Range("K113").Value = 4965.1665
$ws.Range("M113").Value = -2795.1665
$ws.Range("H134").Value = 4720150
$ws.Range("I134").Value = 7143856.5
$ws.Range("K134").Value = 21431569.5
$ws.Range("M134").Value = -21429034.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4208.974
$ws.Range("J31").Value = 7786.5625
$ws.Range("L31").Value = 7786.5625
$ws.Range("N31").Value = -8376.5625
$ws.Range("H34").Value = 4208.974
$ws.Range("J34").Value = 7786.5625
$ws.Range("L34").Value = 7786.5625
$ws.Range("N34").Value = -8190.5625
$ws.Range("H86").Value = 16453619
$ws.Range("I86").Value = 22328238
$ws.Range("K86").Value = 22328238
$ws.Range("M86").Value = -22327115
$ws.Range("H89").Value = 16453619
$ws.Range("I89").Value = 22328238
$ws.Range("K89").Value = 111641190
$ws.Range("M89").Value = -111635574

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1921.1428
$ws.Range("I17").Value = 89.2
$ws.Range("K17").Value = 267.6
$ws.Range("M17").Value = -98.60000000000002
$ws.Range("H26").Value = 614.8
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = 51
$ws.Range("M26").Value = 237
$ws.Range("H32").Value = 142857280
$ws.Range("I32").Value = 333333400
$ws.Range("J32").Value = 199.75
$ws.Range("K32").Value = 1000000200
$ws.Range("L32").Value = 599.25
$ws.Range("M32").Value = -999999917
$ws.Range("N32").Value = -1165.25
$ws.Range("H46").Value = 253748.5
$ws.Range("I46").Value = 500999.5
$ws.Range("K46").Value = 1502998.5
$ws.Range("M46").Value = -1502907.5
$ws.Range("H60").Value = 884.2857
$ws.Range("I60").Value = 875
$ws.Range("J60").Value = 896.6667
$ws.Range("K60").Value = 2625
$ws.Range("L60").Value = 2690.0001
$ws.Range("M60").Value = -2374
$ws.Range("N60").Value = -3192.0001
$ws.Range("H70").Value = 996
$ws.Range("I70").Value = 996
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2988
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2673
$ws.Range("H73").Value = 996
$ws.Range("I73").Value = 996
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2988
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1896
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H75").Value = 35105976
$ws.Range("I75").Value = 37038356
$ws.Range("J75").Value = 33366836
$ws.Range("K75").Value = 111115068
$ws.Range("L75").Value = 100100508
$ws.Range("M75").Value = -111114070
$ws.Range("N75").Value = -100102504
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H78").Value = 35105976
$ws.Range("I78").Value = 37038356
$ws.Range("J78").Value = 33366836
$ws.Range("K78").Value = 333345204
$ws.Range("L78").Value = 300301524
$ws.Range("M78").Value = -333340212
$ws.Range("N78").Value = -300311508
$ws.Range("H122").Value = 1287271.6
$ws.Range("I122").Value = 2572786
$ws.Range("J122").Value = 1757.2727
$ws.Range("K122").Value = 23155074
$ws.Range("L122").Value = 15815.4543
$ws.Range("M122").Value = -23152624
$ws.Range("N122").Value = -20715.4543
$ws.Range("H140").Value = 93610.37
$ws.Range("I140").Value = 168249.33
$ws.Range("J140").Value = 4043.6
$ws.Range("K140").Value = 504747.99
$ws.Range("L140").Value = 12130.8
$ws.Range("M140").Value = -499567.99
$ws.Range("N140").Value = -22490.8
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2942.6843
$ws.Range("I80").Value = 2461.875
$ws.Range("J80").Value = 3292.3635
$ws.Range("K80").Value = 2461.875
$ws.Range("L80").Value = 3292.3635
$ws.Range("M80").Value = -1463.875
$ws.Range("N80").Value = -5288.363499999999
$ws.Range("H83").Value = 2942.6843
$ws.Range("I83").Value = 2461.875
$ws.Range("J83").Value = 3292.3635
$ws.Range("K83").Value = 12309.375
$ws.Range("L83").Value = 16461.8175
$ws.Range("M83").Value = -7317.375
$ws.Range("N83").Value = -26445.8175
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
$ws.Range("H136").Value = 43751.207
$ws.Range("J136").Value = 46501.45
$ws.Range("L136").Value = 139504.35
$ws.Range("N136").Value = -144604.35
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 74732.336
$ws.Range("I74").Value = 74732.336
$ws.Range("K74").Value = 74732.336
$ws.Range("M74").Value = -73734.336
$ws.Range("H77").Value = 74732.336
$ws.Range("I77").Value = 74732.336
$ws.Range("K77").Value = 224197.008
$ws.Range("M77").Value = -219205.008
$ws.Range("H122").Value = 3985.4783
$ws.Range("I122").Value = 2817.7334
$ws.Range("K122").Value = 8453.200199999999
$ws.Range("M122").Value = -6003.200199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28434.273
$ws.Range("I62").Value = 32136.158
$ws.Range("K62").Value = 32136.158
$ws.Range("M62").Value = -31512.158
$ws.Range("H65").Value = 28434.273
$ws.Range("I65").Value = 32136.158
$ws.Range("K65").Value = 160680.79
$ws.Range("M65").Value = -157560.79
$ws.Range("H70").Value = 96018.8
$ws.Range("J70").Value = 112499.75
$ws.Range("L70").Value = 112499.75
$ws.Range("N70").Value = -113129.75
$ws.Range("H73").Value = 96018.8
$ws.Range("J73").Value = 112499.75
$ws.Range("L73").Value = 112499.75
$ws.Range("N73").Value = -114683.75
$ws.Range("H132").Value = 4995.0225
$ws.Range("I132").Value = 4584.5757
$ws.Range("J132").Value = 6123.75
$ws.Range("K132").Value = 13753.7271
$ws.Range("L132").Value = 18371.25
$ws.Range("M132").Value = -11223.7271
$ws.Range("N132").Value = -23431.25
$ws.Range("H136").Value = 12471854
$ws.Range("J136").Value = 440496.53
$ws.Range("L136").Value = 1321489.59
$ws.Range("N136").Value = -1326589.59
